$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -21.718
    13 = -22.219
    16 = -22.027
    18 = -22.128
    20 = -20.385
    26 = -21.284
    27 = -21.91999999999999
    29 = -21.301
    35 = -19.823
    36 = -20.391
    45 = -21.573
    55 = -22.187
    57 = -22.258
    69 = -21.69
    76 = -20.009
    78 = -19.854
    82 = -22.018
    83 = -22.105
    93 = -21.433
    97 = -22.136
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
